$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must stay as literal text (would otherwise be auto-converted
# to number/date/boolean by Excel's smart entry) need NumberFormat "@" first.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("AA4").NumberFormat = "@"

# ---- Row 3 ----
$ws.Range("A3").Value = 112550147
$ws.Range("B3").Value = 78713
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("I3").Value = "1"
$ws.Range("P3").Value = "Flotjärnen, Jmt"
$ws.Range("Q3").Value = 514206
$ws.Range("R3").Value = 7012178
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Ragunda"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Stugun"
$ws.Range("Y3").Value = "2023-10-04"
$ws.Range("AA3").Value = "2023-10-04"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AN3").Value = 1
$ws.Range("AO3").Value = "1 substratenheter"
$ws.Range("AW3").Value = "Sebastian Acker"
$ws.Range("AX3").Value = "Sebastian Acker"
$ws.Range("AY3").Value = "SCA Skog Naturvärdesinventering"

# ---- Row 4 ----
$ws.Range("A4").Value = 112550148
$ws.Range("B4").Value = 78714
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2081
$ws.Range("F4").Value = "Skrovellav"
$ws.Range("G4").Value = "Lobaria scrobiculata"
$ws.Range("H4").Value = "(Scop.) DC."
$ws.Range("I4").Value = "1"
$ws.Range("P4").Value = "Flotjärnen, Jmt"
$ws.Range("Q4").Value = 514204
$ws.Range("R4").Value = 7012186
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Ragunda"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Stugun"
$ws.Range("Y4").Value = "2023-10-04"
$ws.Range("AA4").Value = "2023-10-04"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AN4").Value = 1
$ws.Range("AO4").Value = "1 substratenheter"
$ws.Range("AW4").Value = "Sebastian Acker"
$ws.Range("AX4").Value = "Sebastian Acker"
$ws.Range("AY4").Value = "SCA Skog Naturvärdesinventering"

# AT3 / AT4 are empty placeholder cells (inlineStr with no text) in the
# source data. Give them the text number format so that, if the engine
# materializes a cell for them, it is treated as text rather than a blank
# numeric cell.
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT4").NumberFormat = "@"
